# Agregado de retrazo forzado para mejorar output
#
# Updates the "26-05-2012" performance sheet: rewrites several measured
# input/output throughput samples (rows 7-14 and 18-25), renames the
# "[1000T->P]" scenario label to "[2000T->P]" (rows 14, 25, 36, 47), and
# moves the sheet's top-left/selection back to F7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("26-05-2012")

# --- Block 1 (rows 7-14): Input/Output(tareas/ms) measurements ---
$ws.Range("E7").Value  = 10016
$ws.Range("F7").Value  = 526.59325079872201
$ws.Range("G7").Value  = 526.59305111821004

$ws.Range("F8").Value  = 480.5324
$ws.Range("G8").Value  = 480.53219999999999

$ws.Range("F9").Value  = 481.23910000000001
$ws.Range("G9").Value  = 481.2389

$ws.Range("E10").Value = 10016
$ws.Range("F10").Value = 482.11711261980798
$ws.Range("G10").Value = 482.11691293929698

$ws.Range("F11").Value = 469.99029999999999
$ws.Range("G11").Value = 469.99020000000002

$ws.Range("F12").Value = 465.40199700449301
$ws.Range("G12").Value = 465.40179730404299

$ws.Range("F13").Value = 476.159948243256
$ws.Range("G13").Value = 476.15865432467399

$ws.Range("D14").Value = "[2000T->P]"
$ws.Range("E14").Value = 10547
$ws.Range("F14").Value = 465.42438608135001
$ws.Range("G14").Value = 465.42410164027598

# --- Block 2 (rows 18-25): Input/Output(tareas/ms) measurements ---
$ws.Range("E18").Value = 10000
$ws.Range("F18").Value = 153.9769
$ws.Range("G18").Value = 153.96129999999999

$ws.Range("E19").Value = 10000
$ws.Range("F19").Value = 157.39109999999999
$ws.Range("G19").Value = 157.38749999999999

$ws.Range("E20").Value = 10000
$ws.Range("F20").Value = 157.43100000000001
$ws.Range("G20").Value = 157.42910000000001

$ws.Range("E21").Value = 10031
$ws.Range("F21").Value = 154.53215033396401
$ws.Range("G21").Value = 154.501844282723

$ws.Range("E22").Value = 10016
$ws.Range("F22").Value = 156.50938498402499
$ws.Range("G22").Value = 156.47893370606999

$ws.Range("E23").Value = 10093
$ws.Range("F23").Value = 152.08391954820101
$ws.Range("G23").Value = 152.053799663132

$ws.Range("E24").Value = 10563
$ws.Range("F24").Value = 153.37707090788601
$ws.Range("G24").Value = 153.348196535075

$ws.Range("D25").Value = "[2000T->P]"
$ws.Range("E25").Value = 16391
$ws.Range("F25").Value = 144.614666585321
$ws.Range("G25").Value = 144.59978036727401

# --- Scenario-label rename elsewhere on the sheet ---
$ws.Range("D36").Value = "[2000T->P]"
$ws.Range("D47").Value = "[2000T->P]"

# --- Restore the view: drop the old scrolled position and select F7 ---
$ws.Activate()
$ws.Range("F7").Select()
